{"js": "// The deployment guide gained two new \"(Optional)\" steps for starting /\n// stopping the nginx reverse-proxy helper scripts. Each is inserted as a\n// sub-bullet (same list/level as its sibling steps) right after the\n// matching \"common-data-model\" stop/run step, with the word \"Optional\"\n// in bold inside the trailing parenthetical.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the two anchor paragraphs by their exact text. \"stop.sh\" text is\n// unique, but \"sh ./data-ingestion-suite/docker/deployment/run.sh\" appears\n// twice in the document (once in the earlier onFHIR deployment section and\n// once \u2014 the one we want \u2014 as the final step of \"Run all containers:\" at\n// the very end). Keep scanning forward so we end up with the LAST match.\nlet stopAnchor = null;\nlet runAnchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"sh ./common-data-model/docker/stop.sh\") {\n    stopAnchor = paragraphs.items[i];\n  }\n  if (t === \"sh ./data-ingestion-suite/docker/deployment/run.sh\") {\n    runAnchor = paragraphs.items[i];\n  }\n}\n\nif (!stopAnchor || !runAnchor) {\n  throw new Error(\"Could not locate anchor paragraphs for the edit.\");\n}\n\n// 1) New step right after \"sh ./common-data-model/docker/stop.sh\":\n//    \"sh ./data-ingestion-suite/docker/proxy/stop.sh (Optional)\"\nconst stopOptionalPara = stopAnchor.insertParagraph(\n  \"sh ./data-ingestion-suite/docker/proxy/stop.sh (Optional)\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst stopBoldResults = stopOptionalPara.search(\"Optional\", { matchCase: true });\nstopBoldResults.load(\"text\");\nawait context.sync();\nstopBoldResults.items[0].font.bold = true;\nawait context.sync();\n\n// 2) New step right after the final \"sh .../deployment/run.sh\" step:\n//    \"sh ./data-ingestion-suite/docker/proxy/run.sh (Optional)\"\nconst runOptionalPara = runAnchor.insertParagraph(\n  \"sh ./data-ingestion-suite/docker/proxy/run.sh (Optional)\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst runBoldResults = runOptionalPara.search(\"Optional\", { matchCase: true });\nrunBoldResults.load(\"text\");\nawait context.sync();\nrunBoldResults.items[0].font.bold = true;\nawait context.sync();\n", "ps1": "# The deployment guide gained two new \"(Optional)\" steps for starting /\n# stopping the nginx reverse-proxy helper scripts. Each is inserted as a\n# sub-bullet (same list/level as its sibling steps) right after the\n# matching \"common-data-model\" stop/run step, with the word \"Optional\"\n# in bold inside the trailing parenthetical.\n\n$d = $word.ActiveDocument\n\nfunction Find-LastParagraphIndexByText($doc, $targetText) {\n    $paras = $doc.Paragraphs\n    $result = $null\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $t = $paras.Item($i).Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $targetText) {\n            # Some anchor texts (e.g. the deployment run.sh step) occur more\n            # than once in the document; keep scanning so we land on the\n            # LAST match, which is the one at the end of the \"Run all\n            # containers:\" / \"Stop all containers:\" sections.\n            $result = $i\n        }\n    }\n    return $result\n}\n\nfunction Insert-OptionalStep($doc, $anchorText, $newText) {\n    $anchorIndex = Find-LastParagraphIndexByText $doc $anchorText\n    if ($anchorIndex -eq $null) {\n        throw \"Anchor paragraph not found: $anchorText\"\n    }\n\n    $paras = $doc.Paragraphs\n    $anchor = $paras.Item($anchorIndex)\n    # Inserts a new (initially empty) paragraph right after $anchor, which\n    # inherits $anchor's paragraph formatting (ListParagraph style, same\n    # list level/numId) automatically.\n    $anchor.Range.InsertParagraphAfter()\n\n    $paras = $doc.Paragraphs\n    $newPara = $paras.Item($anchorIndex + 1)\n    $newPara.Range.Text = $newText\n\n    # Bold just the word \"Optional\" inside the new paragraph.\n    $findRange = $newPara.Range.Duplicate()\n    $findRange.Find.Execute(\"Optional\")\n    $findRange.Font.Bold = 1\n}\n\nInsert-OptionalStep $d \"sh ./common-data-model/docker/stop.sh\" \"sh ./data-ingestion-suite/docker/proxy/stop.sh (Optional)\"\nInsert-OptionalStep $d \"sh ./data-ingestion-suite/docker/deployment/run.sh\" \"sh ./data-ingestion-suite/docker/proxy/run.sh (Optional)\"\n"}
